$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45736
$ws.Cells.Item(2, 2).Value = 660.2249999999999
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 5).Value = 176
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 1309.1
$ws.Cells.Item(2, 8).Value = 39
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 595
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 1294.225
$ws.Cells.Item(2, 13).Value = 1485.1
$ws.Cells.Item(2, 14).Value = -190.875
$ws.Cells.Item(3, 1).Value = 45736.01041666666
$ws.Cells.Item(3, 2).Value = 660.2249999999999
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 5).Value = 176
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 1262
$ws.Cells.Item(3, 8).Value = 39
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 595
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 1294.225
$ws.Cells.Item(3, 13).Value = 1438
$ws.Cells.Item(3, 14).Value = -143.7750000000001
$ws.Cells.Item(4, 1).Value = 45736.02083333334
$ws.Cells.Item(4, 2).Value = 660.2249999999999
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 5).Value = 176
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 1251.5
$ws.Cells.Item(4, 8).Value = 39
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 595
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 1294.225
$ws.Cells.Item(4, 13).Value = 1427.5
$ws.Cells.Item(4, 14).Value = -133.2750000000001
$ws.Cells.Item(5, 1).Value = 45736.03125
$ws.Cells.Item(5, 2).Value = 660.2249999999999
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 5).Value = 176
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 1220.8
$ws.Cells.Item(5, 8).Value = 39
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 595
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 1294.225
$ws.Cells.Item(5, 13).Value = 1396.8
$ws.Cells.Item(5, 14).Value = -102.575
$ws.Cells.Item(6, 1).Value = 45736.04166666666
$ws.Cells.Item(6, 2).Value = 612.0500000000002
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 5).Value = 335
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 896.7
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 562
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 1174.05
$ws.Cells.Item(6, 13).Value = 1231.7
$ws.Cells.Item(6, 14).Value = -57.64999999999986
$ws.Cells.Item(7, 1).Value = 45736.05208333334
$ws.Cells.Item(7, 2).Value = 612.0500000000002
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 5).Value = 335
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 921.8
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 562
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 1174.05
$ws.Cells.Item(7, 13).Value = 1256.8
$ws.Cells.Item(7, 14).Value = -82.74999999999977
$ws.Cells.Item(8, 1).Value = 45736.0625
$ws.Cells.Item(8, 2).Value = 612.0500000000002
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 5).Value = 335
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 891.5999999999999
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 562
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 1174.05
$ws.Cells.Item(8, 13).Value = 1226.6
$ws.Cells.Item(8, 14).Value = -52.54999999999973
$ws.Cells.Item(9, 1).Value = 45736.07291666666
$ws.Cells.Item(9, 2).Value = 612.0500000000002
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 5).Value = 335
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 902.7
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 562
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 1174.05
$ws.Cells.Item(9, 13).Value = 1237.7
$ws.Cells.Item(9, 14).Value = -63.64999999999986
$ws.Cells.Item(10, 1).Value = 45736.08333333334
$ws.Cells.Item(10, 2).Value = 625.45
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 5).Value = 357
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 759.8000000000002
$ws.Cells.Item(10, 9).Value = 42
$ws.Cells.Item(10, 10).Value = 546
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 1171.45
$ws.Cells.Item(10, 13).Value = 1158.8
$ws.Cells.Item(10, 14).Value = 12.64999999999986
$ws.Cells.Item(11, 1).Value = 45736.09375
$ws.Cells.Item(11, 2).Value = 625.45
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 5).Value = 357
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 755.3000000000002
$ws.Cells.Item(11, 9).Value = 42
$ws.Cells.Item(11, 10).Value = 546
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 1171.45
$ws.Cells.Item(11, 13).Value = 1154.3
$ws.Cells.Item(11, 14).Value = 17.14999999999986
$ws.Cells.Item(12, 1).Value = 45736.10416666666
$ws.Cells.Item(12, 2).Value = 625.45
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 5).Value = 357
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 786.3000000000002
$ws.Cells.Item(12, 9).Value = 42
$ws.Cells.Item(12, 10).Value = 546
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 1171.45
$ws.Cells.Item(12, 13).Value = 1185.3
$ws.Cells.Item(12, 14).Value = -13.85000000000014
$ws.Cells.Item(13, 1).Value = 45736.11458333334
$ws.Cells.Item(13, 2).Value = 625.45
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 5).Value = 357
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 781
$ws.Cells.Item(13, 9).Value = 42
$ws.Cells.Item(13, 10).Value = 546
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 1171.45
$ws.Cells.Item(13, 13).Value = 1180
$ws.Cells.Item(13, 14).Value = -8.549999999999955
$ws.Cells.Item(14, 1).Value = 45736.125
$ws.Cells.Item(14, 2).Value = 709.875
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 5).Value = 355
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 900.9000000000001
$ws.Cells.Item(14, 9).Value = 58
$ws.Cells.Item(14, 10).Value = 545
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 1254.875
$ws.Cells.Item(14, 13).Value = 1313.9
$ws.Cells.Item(14, 14).Value = -59.02500000000009
$ws.Cells.Item(15, 1).Value = 45736.13541666666
$ws.Cells.Item(15, 2).Value = 709.875
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 5).Value = 355
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 914.4000000000001
$ws.Cells.Item(15, 9).Value = 58
$ws.Cells.Item(15, 10).Value = 545
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 1254.875
$ws.Cells.Item(15, 13).Value = 1327.4
$ws.Cells.Item(15, 14).Value = -72.52500000000009
$ws.Cells.Item(16, 1).Value = 45736.14583333334
$ws.Cells.Item(16, 2).Value = 709.875
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 5).Value = 355
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 944.4000000000001
$ws.Cells.Item(16, 9).Value = 58
$ws.Cells.Item(16, 10).Value = 545
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 1254.875
$ws.Cells.Item(16, 13).Value = 1357.4
$ws.Cells.Item(16, 14).Value = -102.5250000000001
$ws.Cells.Item(17, 1).Value = 45736.15625
$ws.Cells.Item(17, 2).Value = 709.875
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 5).Value = 355
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 925.1000000000001
$ws.Cells.Item(17, 9).Value = 58
$ws.Cells.Item(17, 10).Value = 545
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1254.875
$ws.Cells.Item(17, 13).Value = 1338.1
$ws.Cells.Item(17, 14).Value = -83.22500000000014
$ws.Cells.Item(18, 1).Value = 45736.16666666666
$ws.Cells.Item(18, 2).Value = 848.9000000000001
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 5).Value = 340
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 1150.5
$ws.Cells.Item(18, 8).Value = 6
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 594
$ws.Cells.Item(18, 12).Value = 1448.9
$ws.Cells.Item(18, 13).Value = 1490.5
$ws.Cells.Item(18, 14).Value = -41.59999999999991
$ws.Cells.Item(19, 1).Value = 45736.17708333334
$ws.Cells.Item(19, 2).Value = 848.9000000000001
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 5).Value = 340
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 1158.7
$ws.Cells.Item(19, 8).Value = 6
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 594
$ws.Cells.Item(19, 12).Value = 1448.9
$ws.Cells.Item(19, 13).Value = 1498.7
$ws.Cells.Item(19, 14).Value = -49.79999999999995
$ws.Cells.Item(20, 1).Value = 45736.1875
$ws.Cells.Item(20, 2).Value = 848.9000000000001
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 5).Value = 340
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 1242.2
$ws.Cells.Item(20, 8).Value = 6
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 594
$ws.Cells.Item(20, 12).Value = 1448.9
$ws.Cells.Item(20, 13).Value = 1582.2
$ws.Cells.Item(20, 14).Value = -133.3
$ws.Cells.Item(21, 1).Value = 45736.19791666666
$ws.Cells.Item(21, 2).Value = 848.9000000000001
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 5).Value = 340
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 1262.6
$ws.Cells.Item(21, 8).Value = 6
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 594
$ws.Cells.Item(21, 12).Value = 1448.9
$ws.Cells.Item(21, 13).Value = 1602.6
$ws.Cells.Item(21, 14).Value = -153.6999999999998
$ws.Cells.Item(22, 1).Value = 45736.20833333334
$ws.Cells.Item(22, 2).Value = 631.5250000000001
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 54
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 1246.7
$ws.Cells.Item(22, 8).Value = 77
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 547
$ws.Cells.Item(22, 12).Value = 1255.525
$ws.Cells.Item(22, 13).Value = 1300.7
$ws.Cells.Item(22, 14).Value = -45.17499999999973
$ws.Cells.Item(23, 1).Value = 45736.21875
$ws.Cells.Item(23, 2).Value = 631.5250000000001
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 54
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 1307
$ws.Cells.Item(23, 8).Value = 77
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 547
$ws.Cells.Item(23, 12).Value = 1255.525
$ws.Cells.Item(23, 13).Value = 1361
$ws.Cells.Item(23, 14).Value = -105.4749999999999
$ws.Cells.Item(24, 1).Value = 45736.22916666666
$ws.Cells.Item(24, 2).Value = 631.5250000000001
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 54
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 1409.2
$ws.Cells.Item(24, 8).Value = 77
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 547
$ws.Cells.Item(24, 12).Value = 1255.525
$ws.Cells.Item(24, 13).Value = 1463.2
$ws.Cells.Item(24, 14).Value = -207.6749999999997
$ws.Cells.Item(25, 1).Value = 45736.23958333334
$ws.Cells.Item(25, 2).Value = 631.5250000000001
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 54
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 1449.7
$ws.Cells.Item(25, 8).Value = 77
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 547
$ws.Cells.Item(25, 12).Value = 1255.525
$ws.Cells.Item(25, 13).Value = 1503.7
$ws.Cells.Item(25, 14).Value = -248.1749999999997
$ws.Cells.Item(26, 1).Value = 45736.25
$ws.Cells.Item(26, 2).Value = 246.95
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 113
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 648.5999999999999
$ws.Cells.Item(26, 9).Value = 147
$ws.Cells.Item(26, 10).Value = 268
$ws.Cells.Item(26, 12).Value = 627.95
$ws.Cells.Item(26, 13).Value = 795.5999999999999
$ws.Cells.Item(26, 14).Value = -167.6499999999999
$ws.Cells.Item(27, 1).Value = 45736.26041666666
$ws.Cells.Item(27, 2).Value = 246.95
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 113
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 723.3
$ws.Cells.Item(27, 9).Value = 147
$ws.Cells.Item(27, 10).Value = 268
$ws.Cells.Item(27, 12).Value = 627.95
$ws.Cells.Item(27, 13).Value = 870.3
$ws.Cells.Item(27, 14).Value = -242.3499999999999
$ws.Cells.Item(28, 1).Value = 45736.27083333334
$ws.Cells.Item(28, 2).Value = 246.95
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 113
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 550.1999999999999
$ws.Cells.Item(28, 9).Value = 147
$ws.Cells.Item(28, 10).Value = 268
$ws.Cells.Item(28, 12).Value = 627.95
$ws.Cells.Item(28, 13).Value = 697.1999999999999
$ws.Cells.Item(28, 14).Value = -69.24999999999989
$ws.Cells.Item(29, 1).Value = 45736.28125
$ws.Cells.Item(29, 2).Value = 246.95
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 113
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 508.3
$ws.Cells.Item(29, 9).Value = 147
$ws.Cells.Item(29, 10).Value = 268
$ws.Cells.Item(29, 12).Value = 627.95
$ws.Cells.Item(29, 13).Value = 655.3
$ws.Cells.Item(29, 14).Value = -27.34999999999991
$ws.Cells.Item(30, 1).Value = 45736.29166666666
$ws.Cells.Item(30, 3).Value = 323.575
$ws.Cells.Item(30, 4).Value = 134
$ws.Cells.Item(30, 6).Value = 373.8
$ws.Cells.Item(30, 9).Value = 331
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 12
$ws.Cells.Item(30, 12).Value = 507.8
$ws.Cells.Item(30, 13).Value = 666.575
$ws.Cells.Item(30, 14).Value = -158.775
$ws.Cells.Item(31, 1).Value = 45736.30208333334
$ws.Cells.Item(31, 3).Value = 323.575
$ws.Cells.Item(31, 4).Value = 134
$ws.Cells.Item(31, 6).Value = 392.2
$ws.Cells.Item(31, 9).Value = 331
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 12
$ws.Cells.Item(31, 12).Value = 526.2
$ws.Cells.Item(31, 13).Value = 666.575
$ws.Cells.Item(31, 14).Value = -140.375
$ws.Cells.Item(32, 1).Value = 45736.3125
$ws.Cells.Item(32, 3).Value = 323.575
$ws.Cells.Item(32, 4).Value = 134
$ws.Cells.Item(32, 6).Value = 645.1
$ws.Cells.Item(32, 9).Value = 331
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 12
$ws.Cells.Item(32, 12).Value = 779.1
$ws.Cells.Item(32, 13).Value = 666.575
$ws.Cells.Item(32, 14).Value = 112.525
$ws.Cells.Item(33, 1).Value = 45736.32291666666
$ws.Cells.Item(33, 3).Value = 323.575
$ws.Cells.Item(33, 4).Value = 134
$ws.Cells.Item(33, 6).Value = 843.2
$ws.Cells.Item(33, 9).Value = 331
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 12
$ws.Cells.Item(33, 12).Value = 977.2
$ws.Cells.Item(33, 13).Value = 666.575
$ws.Cells.Item(33, 14).Value = 310.625
$ws.Cells.Item(34, 1).Value = 45736.33333333334
$ws.Cells.Item(34, 3).Value = 197.2249999999999
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(34, 6).Value = 653
$ws.Cells.Item(34, 9).Value = 269
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 116
$ws.Cells.Item(34, 12).Value = 654
$ws.Cells.Item(34, 13).Value = 582.2249999999999
$ws.Cells.Item(34, 14).Value = 71.77500000000009
$ws.Cells.Item(35, 1).Value = 45736.34375
$ws.Cells.Item(35, 3).Value = 197.2249999999999
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 6).Value = 592
$ws.Cells.Item(35, 9).Value = 269
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 116
$ws.Cells.Item(35, 12).Value = 593
$ws.Cells.Item(35, 13).Value = 582.2249999999999
$ws.Cells.Item(35, 14).Value = 10.77500000000009
$ws.Cells.Item(36, 1).Value = 45736.35416666666
$ws.Cells.Item(36, 3).Value = 197.2249999999999
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 6).Value = 643.2
$ws.Cells.Item(36, 9).Value = 269
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 116
$ws.Cells.Item(36, 12).Value = 644.2
$ws.Cells.Item(36, 13).Value = 582.2249999999999
$ws.Cells.Item(36, 14).Value = 61.97500000000014
$ws.Cells.Item(37, 1).Value = 45736.36458333334
$ws.Cells.Item(37, 3).Value = 197.2249999999999
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(37, 6).Value = 773.6
$ws.Cells.Item(37, 9).Value = 269
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 116
$ws.Cells.Item(37, 12).Value = 774.6
$ws.Cells.Item(37, 13).Value = 582.2249999999999
$ws.Cells.Item(37, 14).Value = 192.3750000000001
$ws.Cells.Item(38, 1).Value = 45736.375
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(38, 3).Value = 652.2
$ws.Cells.Item(38, 4).Value = 122
$ws.Cells.Item(38, 6).Value = 1418
$ws.Cells.Item(38, 9).Value = 350
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 263
$ws.Cells.Item(38, 12).Value = 1540
$ws.Cells.Item(38, 13).Value = 1265.2
$ws.Cells.Item(38, 14).Value = 274.8
$ws.Cells.Item(39, 1).Value = 45736.38541666666
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(39, 3).Value = 652.2
$ws.Cells.Item(39, 4).Value = 122
$ws.Cells.Item(39, 6).Value = 1471.8
$ws.Cells.Item(39, 9).Value = 350
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 263
$ws.Cells.Item(39, 12).Value = 1593.8
$ws.Cells.Item(39, 13).Value = 1265.2
$ws.Cells.Item(39, 14).Value = 328.5999999999999
$ws.Cells.Item(40, 1).Value = 45736.39583333334
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(40, 3).Value = 652.2
$ws.Cells.Item(40, 4).Value = 122
$ws.Cells.Item(40, 6).Value = 1443.2
$ws.Cells.Item(40, 9).Value = 350
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 263
$ws.Cells.Item(40, 12).Value = 1565.2
$ws.Cells.Item(40, 13).Value = 1265.2
$ws.Cells.Item(40, 14).Value = 300
$ws.Cells.Item(41, 1).Value = 45736.40625
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 652.2
$ws.Cells.Item(41, 4).Value = 122
$ws.Cells.Item(41, 6).Value = 1506.4
$ws.Cells.Item(41, 9).Value = 350
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 263
$ws.Cells.Item(41, 12).Value = 1628.4
$ws.Cells.Item(41, 13).Value = 1265.2
$ws.Cells.Item(41, 14).Value = 363.2
$ws.Cells.Item(42, 1).Value = 45736.41666666666
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 480.05
$ws.Cells.Item(42, 4).Value = 190
$ws.Cells.Item(42, 6).Value = 1201
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 9).Value = 253
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 172
$ws.Cells.Item(42, 12).Value = 1391
$ws.Cells.Item(42, 13).Value = 905.05
$ws.Cells.Item(42, 14).Value = 485.95
$ws.Cells.Item(43, 1).Value = 45736.42708333334
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(43, 3).Value = 480.05
$ws.Cells.Item(43, 4).Value = 190
$ws.Cells.Item(43, 6).Value = 1196
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 9).Value = 253
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 172
$ws.Cells.Item(43, 12).Value = 1386
$ws.Cells.Item(43, 13).Value = 905.05
$ws.Cells.Item(43, 14).Value = 480.95
$ws.Cells.Item(44, 1).Value = 45736.4375
$ws.Cells.Item(44, 2).Value = 0
$ws.Cells.Item(44, 3).Value = 480.05
$ws.Cells.Item(44, 4).Value = 190
$ws.Cells.Item(44, 6).Value = 1208
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 9).Value = 253
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 172
$ws.Cells.Item(44, 12).Value = 1398
$ws.Cells.Item(44, 13).Value = 905.05
$ws.Cells.Item(44, 14).Value = 492.95
$ws.Cells.Item(45, 1).Value = 45736.44791666666
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = 480.05
$ws.Cells.Item(45, 4).Value = 190
$ws.Cells.Item(45, 6).Value = 1211
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 9).Value = 253
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 172
$ws.Cells.Item(45, 12).Value = 1401
$ws.Cells.Item(45, 13).Value = 905.05
$ws.Cells.Item(45, 14).Value = 495.95
$ws.Cells.Item(46, 1).Value = 45736.45833333334
$ws.Cells.Item(46, 3).Value = 159.125
$ws.Cells.Item(46, 4).Value = 192
$ws.Cells.Item(46, 6).Value = 612.3
$ws.Cells.Item(46, 9).Value = 99
$ws.Cells.Item(46, 11).Value = 62
$ws.Cells.Item(46, 12).Value = 804.3
$ws.Cells.Item(46, 13).Value = 320.125
$ws.Cells.Item(46, 14).Value = 484.175
$ws.Cells.Item(47, 1).Value = 45736.46875
$ws.Cells.Item(47, 3).Value = 159.125
$ws.Cells.Item(47, 4).Value = 192
$ws.Cells.Item(47, 6).Value = 582.0999999999999
$ws.Cells.Item(47, 9).Value = 99
$ws.Cells.Item(47, 11).Value = 62
$ws.Cells.Item(47, 12).Value = 774.0999999999999
$ws.Cells.Item(47, 13).Value = 320.125
$ws.Cells.Item(47, 14).Value = 453.9749999999999
$ws.Cells.Item(48, 1).Value = 45736.47916666666
$ws.Cells.Item(48, 3).Value = 159.125
$ws.Cells.Item(48, 4).Value = 192
$ws.Cells.Item(48, 6).Value = 535.7
$ws.Cells.Item(48, 9).Value = 99
$ws.Cells.Item(48, 11).Value = 62
$ws.Cells.Item(48, 12).Value = 727.7
$ws.Cells.Item(48, 13).Value = 320.125
$ws.Cells.Item(48, 14).Value = 407.575
$ws.Cells.Item(49, 1).Value = 45736.48958333334
$ws.Cells.Item(49, 3).Value = 159.125
$ws.Cells.Item(49, 4).Value = 192
$ws.Cells.Item(49, 6).Value = 542.7
$ws.Cells.Item(49, 9).Value = 99
$ws.Cells.Item(49, 11).Value = 62
$ws.Cells.Item(49, 12).Value = 734.7
$ws.Cells.Item(49, 13).Value = 320.125
$ws.Cells.Item(49, 14).Value = 414.575
$ws.Cells.Item(50, 1).Value = 45736.5
$ws.Cells.Item(50, 6).Value = 449.5
$ws.Cells.Item(50, 12).Value = 449.5
$ws.Cells.Item(50, 14).Value = 449.5
$ws.Cells.Item(51, 1).Value = 45736.51041666666
$ws.Cells.Item(51, 6).Value = 454.3000000000001
$ws.Cells.Item(51, 12).Value = 454.3000000000001
$ws.Cells.Item(51, 14).Value = 454.3000000000001
$ws.Cells.Item(52, 1).Value = 45736.52083333334
$ws.Cells.Item(53, 1).Value = 45736.53125
$ws.Cells.Item(54, 1).Value = 45736.54166666666
$ws.Cells.Item(55, 1).Value = 45736.55208333334
$ws.Cells.Item(56, 1).Value = 45736.5625
$ws.Cells.Item(57, 1).Value = 45736.57291666666
$ws.Cells.Item(58, 1).Value = 45736.58333333334
$ws.Cells.Item(59, 1).Value = 45736.59375
$ws.Cells.Item(60, 1).Value = 45736.60416666666
$ws.Cells.Item(61, 1).Value = 45736.61458333334
$ws.Cells.Item(62, 1).Value = 45736.625
$ws.Cells.Item(63, 1).Value = 45736.63541666666
$ws.Cells.Item(64, 1).Value = 45736.64583333334
$ws.Cells.Item(65, 1).Value = 45736.65625
$ws.Cells.Item(66, 1).Value = 45736.66666666666
$ws.Cells.Item(67, 1).Value = 45736.67708333334
$ws.Cells.Item(68, 1).Value = 45736.6875
$ws.Cells.Item(69, 1).Value = 45736.69791666666
$ws.Cells.Item(70, 1).Value = 45736.70833333334
$ws.Cells.Item(71, 1).Value = 45736.71875
$ws.Cells.Item(72, 1).Value = 45736.72916666666
$ws.Cells.Item(73, 1).Value = 45736.73958333334
$ws.Cells.Item(74, 1).Value = 45736.75
$ws.Cells.Item(75, 1).Value = 45736.76041666666
$ws.Cells.Item(76, 1).Value = 45736.77083333334
$ws.Cells.Item(77, 1).Value = 45736.78125
$ws.Cells.Item(78, 1).Value = 45736.79166666666
$ws.Cells.Item(79, 1).Value = 45736.80208333334
$ws.Cells.Item(80, 1).Value = 45736.8125
$ws.Cells.Item(81, 1).Value = 45736.82291666666
$ws.Cells.Item(82, 1).Value = 45736.83333333334
$ws.Cells.Item(83, 1).Value = 45736.84375
$ws.Cells.Item(84, 1).Value = 45736.85416666666
$ws.Cells.Item(85, 1).Value = 45736.86458333334
$ws.Cells.Item(86, 1).Value = 45736.875
$ws.Cells.Item(87, 1).Value = 45736.88541666666
$ws.Cells.Item(88, 1).Value = 45736.89583333334
$ws.Cells.Item(89, 1).Value = 45736.90625
$ws.Cells.Item(90, 1).Value = 45736.91666666666
$ws.Cells.Item(91, 1).Value = 45736.92708333334
$ws.Cells.Item(92, 1).Value = 45736.9375
$ws.Cells.Item(93, 1).Value = 45736.94791666666
$ws.Cells.Item(94, 1).Value = 45736.95833333334
$ws.Cells.Item(95, 1).Value = 45736.96875
$ws.Cells.Item(96, 1).Value = 45736.97916666666
$ws.Cells.Item(97, 1).Value = 45736.98958333334
